# Update "paises.xlsx" (country COVID stats) to the next data refresh:
#  - bump the "last updated" timestamp
#  - refresh the handful of countries whose counts actually changed
#  - re-sort the country table by "Casos totales" (column B) descending,
#    since a couple of countries (India, Costa de Marfil) overtook their
#    neighbours and need to move up in the ranking

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh the "Datos actualizados ..." timestamp banner -------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 20:52"

# --- 2. Push new case numbers for the countries that changed ---------------
# Helper to locate a country's row by name in column A, then write its
# Casos totales / Nuevos casos / Casos activos / Recuperados /
# Casos criticos / Muertes hoy / Muertes.
function Set-CountryRow($country, $b, $c, $d, $e, $f, $g, $h) {
    $found = $ws.Columns("A").Find($country)
    $r = $found.Row
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

Set-CountryRow "Estados Unidos"   603496 16555 38144 540157 12828 1555 25195
Set-CountryRow "India"            11487  1034  1359  9735      0   35   393
Set-CountryRow "Peru"             10303   519  2642  7431    143   14   230
Set-CountryRow "Chequia"           6111    52   642  5308     92   18   161
Set-CountryRow "Lituania"          1070     8   101   940     14    5    29
Set-CountryRow "Tunez"              747    21    43   670     89    0    34
Set-CountryRow "Costa de Marfil"    638    12   114   518      0    0     6
Set-CountryRow "Montenegro"         283     9    46   233      7    1     4
Set-CountryRow "Maldivas"            20     0    16     4      0    0     0
Set-CountryRow "Angola"              19     0     5    12      0    0     2

# --- 3. Re-sort the table by Casos totales (column B), descending ---------
$dataRange = $ws.Range("A4:H216")
$dataRange.Sort($ws.Range("B4"), 2)
